$d = $word.ActiveDocument

# --- Step 1: handle the multi-run "Avaliacao" list paragraph (Metodo/Criterio/
# Norma de recuperacao) first, via precise character sub-ranges, so its internal
# bold labels and w:br line breaks stay untouched. We process the three value
# runs from the END of the paragraph backwards so that earlier offsets in the
# same paragraph remain valid after each edit (lengths change).
$avaliacaoPara = $d.Paragraphs.Item(17)
$baseStart = $avaliacaoPara.Range.Start

# Relative offsets (within this paragraph) of the 3 "value" runs in the original doc:
#   value1 (after "Metodo: ")               -> [8, 82)
#   value2 (after "Criterio: ")              -> [92, 216)
#   value3 (after "Norma de recuperacao: ")  -> [238, 374)

$value3Range = $d.Range($baseStart + 238, $baseStart + 374)
$value3Range.Text = '1) CREMASCO, M. A. Fundamentos de Transferência de Massa, 3ª ed. São Paulo: Editora Blucher; 2021. 2) INCROPERA, F. P.; WITT, D. P. Fundamentos de Transferência de Calor e Massa. 8ª ed. Rio de Janeiro: LTC, 2019. 3) Bird, R. B.; STEWART, W. E.; LIGHTFOOT, E. N. Fenômenos de Transporte. 2ª ed. Rio de Janeiro: LTC, 2004. 4) COULSON, J. M.; RICHARDSON, J. F.; BACKHURST, J. R.; HARKER, J. H. Fluid Flow, Heat Transfer and Mass Transfer. In: COULSON & Richardson Series - Chemical Engineering. 6th ed. Pergamon Press, Oxford, 1999. v.1 5) PERRY''s Chemical Engineers Handbook. Editor in Chief Don W. Green; Late Editor Robert H. Perry. 9ª ed. New York: McGraw-Hill, 2019. 6) WELTY, J. R.; PIGFORD, R. L.; WILKE, C. R. Fundamentals of Momentum, Heat, and Mass Transfer. 6th ed. USA: John Wiley & Sons, Inc, 2014. 7) POLING, B. E.; PRAUSNITZ, J. M.; O''CONNELL, J. The Properties of Gases and Liquids. 5th ed. New York: McGraw-Hill, 2004. 8) CALDAS, J. N.; DE LACERDA, A. I.; VELOSO, E.; PASCHOAL, L. C. M. Internos de Torres: Pratos & Recheios. 2ª ed. Rio de Janeiro: Editora Interciência, 2007.'

$value2Range = $d.Range($baseStart + 92, $baseStart + 216)
$value2Range.Text = ('A recuperação será feita por meio de uma Prova Escrita (PE) e a Média de Recuperação (MR) será calculada pela fórmula: MR = (NF + PE)/2.' + [char]11 + '')

$value1Range = $d.Range($baseStart + 8, $baseStart + 82)
$value1Range.Text = ('A Nota Final (NF) será calculada da seguinte maneira: NF = (P1 + 2*P2)/3' + [char]11 + 'P2 = Nota da Prova (80%) e Nota do Trabalho (20%).' + [char]11 + '')

# --- Step 2: whole-paragraph text swaps for the remaining single-run paragraphs.
# Paragraph identity / Word paragraph-collection index is stable across these
# edits (count of paragraphs never changes), so 1-based Paragraphs.Item(n) indices
# below are safe to use in any order.

# Objetivos (PT) paragraph -> short summary program (PT)
$d.Paragraphs.Item(6).Range.Text = ('1- Introdução:' + [char]11 + '2- Coeficiente de difusão:' + [char]11 + '3- Concentrações, velocidade e fluxos:' + [char]11 + '4 -Equações da continuidade em transferência de massa:' + [char]11 + '5- Difusão em regime permanente sem reação química:' + [char]11 + '6- Difusão com reação química:' + [char]11 + '7- Transferência de massa entre fases.')

# Objetivos (EN) paragraph -> short summary program (EN)
$d.Paragraphs.Item(7).Range.Text = '1 - Introduction: 2 - Diffusion coefficient: 3 - Concentrations, and flow rate: 4 - Equation of continuity for mass transfer: 5 - Diffusion in continuous operation without chemical reaction: 6 - Diffusion with chemical reaction: 7 - Mass transfer between phases.'

# Docente(s) list paragraph -> Objetivos (PT) text
$d.Paragraphs.Item(9).Range.Text = 'Proporcionar ao graduando conhecimentos da teoria básica dos conceitos de transferência de massa com posterior aplicação aos balanços de massa visando obtenção, para os diversos processos físicos e químicos, em particularidade os sistemas estagnados e convectivos, conhecimento do fluxo de transferência de massa, do perfil de concentração, das resistências que prediz o transporte entre as fases.'

# short summary program (PT) -> full Programa text (PT)
$d.Paragraphs.Item(11).Range.Text = '1- Introdução: Transferência de massa: Definição. Classificação das operações que envolvem transferência de massa. Contribuições à transferência de massa. Tipos de difusão. 2- Coeficiente e mecanismos de difusão: Considerações a respeito. Difusão em gases: Análise da primeira lei de Fick. O coeficiente de difusão para gases. Estimativa do coeficiente de difusão a partir de um coeficiente de difusão conhecido em outra temperatura e pressão. Coeficiente de difusão de um soluto em uma mistura gasosa estagnada de multicomponentes. Difusão em líquidos. Difusão em sólidos. 3- Concentrações, velocidades e fluxos: Concentração. Velocidade. Fluxo. A equação de Stefan – Maxwel. Coeficiente convectivo de transferência de massa 4 - Equações da continuidade em transferência de massa: Considerações a respeito. Equações da continuidade mássica e molar de um soluto. Equações da continuidade do soluto A em termos da lei ordinária da difusão. Condições de contorno. 5- Difusão em regime permanente sem reação química: Difusão Unidimensional em regime permanente. Difusão através de filme gasoso inerte e estagnado. Difusão pseudo-estacionária num filme gasoso estagnado. Contradifusão equimolar. Taxa molar em esferas isoladas. Difusão em membranas. 6- Difusão em regime permanente com reação química: Difusão em regime permanente com reação química heterogênea na superfície de uma partícula catalítica não porosa. Difusão com reação química heterogênea na superfície de uma partícula não catalítica e não porosa. Difusão intraparticular com reação química heterogênea. Difusão em regime permanente com reação química homogênea. 7- Transferência de massa entre fases: Considerações a respeito. Técnicas de separação. Transferência de massa entre fases. Teoria das duas resistências. Coeficientes globais de transferência de massa. Coeficientes volumétricos de transferência de massa para torres de recheios. Balanço macroscópio de matéria em equipamentos de separação. Operações contínuas (contracorrente e paralelo). Cálculo da altura efetiva e do diâmetro de uma coluna para operação contínua em um sistema diluído.'

# short summary program (EN) -> Objetivos (EN) text
$d.Paragraphs.Item(12).Range.Text = 'Providing the student knowledge of basic concepts of the theory of mass transfer with subsequent application to obtain mass balances aiming for the different physical and chemical processes, in peculiarity stagnant and convective systems, knowledge of the flow of mass transfer, the profile concentration of resistors that predicts the transport between the phases.'

# full Programa text (PT) -> Metodo value (old)
$d.Paragraphs.Item(14).Range.Text = 'A avaliação será feita por meio de provas escritas e trabalhos em grupos.'

# Bibliografia paragraph -> Docente(s) list text
$d.Paragraphs.Item(19).Range.Text = '5840841 - Gilberto Garcia Cortez'

